# Atualização de bases das ligas, do dia: 29-03-2024 às 17:05
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rows 9 and 10 had their match data (everything except the running
#    index kept in column A) swapped between each other.
$row9  = $ws.Range("B9:AC9").Value2
$row10 = $ws.Range("B10:AC10").Value2

$ws.Range("B9:AC9").Value2  = $row10
$ws.Range("B10:AC10").Value2 = $row9

# 2) The fixture that used to live in row 140 (id 7952735,
#    Zeljeznicar vs Velez Mostar) was removed from the source feed.
#    The remaining fixtures (previously rows 141-143) move up one row,
#    but the running index in column A is NOT touched - only columns
#    B:AA shift - and the now-superfluous last row (143) is removed.
$src = $ws.Range("B141:AA143").Value2
$ws.Range("B140:AA142").Value2 = $src
$ws.Rows("143").Delete() | Out-Null
